$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.514.02"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.913.87"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.70%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.65"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4821"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4069"
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08145"
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.013"
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.43"
$ws.Range("E11").Value = "  +3.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.920.08"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.990"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.118"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.24"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06792"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.68"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.519.14"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.625"
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.77"
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.185"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.094.03"
$ws.Range("E25").Value = "  -1.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.26"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.460"
$ws.Range("E27").Value = "  +6.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.05"
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.100"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.69"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.028"
$ws.Range("E31").Value = "  -4.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09544"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.507"
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.567"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.391"
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06098"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.178"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5938"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.72"
$ws.Range("E40").Value = "  +5.52%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.938"
$ws.Range("E41").Value = "  -4.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1851"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.443"
$ws.Range("E43").Value = "  -4.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.287"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07640"
$ws.Range("E45").Value = "  -3.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.44"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5573"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.939"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "115.91"
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.39"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.051"
$ws.Range("E51").Value = "  +1.62%  "
